$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,@("ECs", "Slit3", "Robo4", "ECs", 3, 1, 3.590118666666667, 10.770356, 0.03935865391742773, 0.03935865391742773, 3, 1, 43.55927533333334, 130.677826, 0.9894183625413969, 0.9894183625413967, 156.3829674806729, 1407.446707326056, 0.03894217491081488, 0.03894217491081488)
    ,@("ECs", "Slit3", "Robo4", "FAPs", 3, 1, 3.590118666666667, 10.770356, 0.03935865391742773, 0.03935865391742773, 2, 0.6666666666666666, 0.205596, 0.616788, 0.004669968820840217, 0.004669968820840216, 0.738114037392, 6.643026336528, 0.0001838036866246282, 0.0001838036866246281)
    ,@("ECs", "Slit3", "Robo4", "sCs", 3, 1, 3.590118666666667, 10.770356, 0.03935865391742773, 0.03935865391742773, 2, 0.6666666666666666, 0.260262, 0.780786, 0.005911668637762975, 0.005911668637762974, 0.934371464424, 8.409343179816, 0.0002326753199882243, 0.0002326753199882243)
    ,@("FAPs", "Slit3", "Robo4", "ECs", 3, 1, 76.92488366666667, 230.774651, 0.8433314202078528, 0.8433314202078527, 3, 1, 43.55927533333334, 130.677826, 0.9894183625413969, 0.9894183625413967, 3350.79218762097, 30157.12968858873, 0.8344075928617645, 0.8344075928617641)
    ,@("FAPs", "Slit3", "Robo4", "FAPs", 3, 1, 76.92488366666667, 230.774651, 0.8433314202078528, 0.8433314202078527, 2, 0.6666666666666666, 0.205596, 0.616788, 0.004669968820840217, 0.004669968820840216, 15.815448382332, 142.339035440988, 0.003938331438005572, 0.003938331438005571)
    ,@("FAPs", "Slit3", "Robo4", "sCs", 3, 1, 76.92488366666667, 230.774651, 0.8433314202078528, 0.8433314202078527, 2, 0.6666666666666666, 0.260262, 0.780786, 0.005911668637762975, 0.005911668637762974, 20.020624072854, 180.185616655686, 0.004985495908082871, 0.004985495908082871)
    ,@("M2", "Slit3", "Robo4", "ECs", 3, 1, 0.041643, 0.124929, 0.0004565343313861054, 0.0004565343313861054, 3, 1, 43.55927533333334, 130.677826, 0.9894183625413969, 0.9894183625413967, 1.813938902706, 16.325450124354, 0.0004517034506039719, 0.0004517034506039718)
    ,@("M2", "Slit3", "Robo4", "FAPs", 3, 1, 0.041643, 0.124929, 0.0004565343313861054, 0.0004565343313861054, 2, 0.6666666666666666, 0.205596, 0.616788, 0.004669968820840217, 0.004669968820840216, 0.008561634228, 0.077054708052, 0.000002132001093216248, 0.000002132001093216247)
    ,@("M2", "Slit3", "Robo4", "sCs", 3, 1, 0.041643, 0.124929, 0.0004565343313861054, 0.0004565343313861054, 2, 0.6666666666666666, 0.260262, 0.780786, 0.005911668637762975, 0.005911668637762974, 0.010838090466, 0.097542814194, 0.000002698879688917328, 0.000002698879688917328)
    ,@("sCs", "Slit3", "Robo4", "ECs", 3, 1, 10.65883866666667, 31.976516, 0.1168533915433334, 0.1168533915433334, 3, 1, 43.55927533333334, 130.677826, 0.9894183625413969, 0.9894183625413967, 464.2912882149129, 4178.621593934216, 0.1156168913182137, 0.1156168913182136)
    ,@("sCs", "Slit3", "Robo4", "FAPs", 3, 1, 10.65883866666667, 31.976516, 0.1168533915433334, 0.1168533915433334, 2, 0.6666666666666666, 0.205596, 0.616788, 0.004669968820840217, 0.004669968820840216, 2.191414594512, 19.722731350608, 0.000545701695116801, 0.0005457016951168008)
    ,@("sCs", "Slit3", "Robo4", "sCs", 3, 1, 10.65883866666667, 31.976516, 0.1168533915433334, 0.1168533915433334, 2, 0.6666666666666666, 0.260262, 0.780786, 0.005911668637762975, 0.005911668637762974, 2.774090669064, 24.966816021576, 0.0006907985300029615, 0.0006907985300029613)
)

$rowIndex = 2
foreach ($row in $data) {
    $ws.Cells.Item($rowIndex, 1).Value = $row[0]
    $ws.Cells.Item($rowIndex, 2).Value = $row[1]
    $ws.Cells.Item($rowIndex, 3).Value = $row[2]
    $ws.Cells.Item($rowIndex, 4).Value = $row[3]
    for ($c = 5; $c -le 20; $c++) {
        $ws.Cells.Item($rowIndex, $c).Value = $row[$c - 1]
    }
    $rowIndex++
}
